$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reqs")
Write-Host $ws.Name
